$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 2 and 3 (Rv1430 <-> Rv2539c) and append new rows 5-13.
$data = @(
    @("Rv2539c", 5, "aroK Rv2539c MTCY159.17", "FUNCTION: Catalyzes the specific phosphorylation of the 3-hydroxyl group of shikimic acid using ATP as a cosubstrate. {ECO:0000269|PubMed:11483005, ECO:0000269|PubMed:17020768}.", 31),
    @("Rv1430", 4, "PE16 Rv1430", "FUNCTION: Esterase that hydrolyzes short to medium chain fatty acid esters with the highest specific activity for p-nitrophenyl caproate (pNPC6). Has lower activity with p-nitrophenyl caprylate (pNPC8) and p-nitrophenyl butyrate (pNPC4). Has weak activity with p-nitrophenyl caprate (pNPC10) and p-nitrophenyl laurate (pNPC12). Does not possess lipolytic activity and cutinase activity. {ECO:0000269|PubMed:23383323}.", 31),
    @("Rv2244", 4, "acpM Rv2244 MTCY427.25", "FUNCTION: Acyl carrier protein involved in meromycolate extension. {ECO:0000250}.", 31),
    @("Rv3923c", 3, "rnpA Rv3923c MTV028.14c", "FUNCTION: RNaseP catalyzes the removal of the 5'-leader sequence from pre-tRNA to produce the mature 5'-terminus. It can also cleave other RNA substrates such as 4.5S RNA. The protein component plays an auxiliary but essential role in vivo by binding to the 5'-leader sequence and broadening the substrate specificity of the ribozyme. {ECO:0000255|HAMAP-Rule:MF_00227}.", 31),
    @("Rv3917c", 3, "parB Rv3917c MTV028.08c", "FUNCTION: Involved in chromosome partition. Localize to both poles of the predivisional cell following completion of DNA replication. Binds to the DNA origin of replication (By similarity). {ECO:0000250}.", 31),
    @("Rv3866", 3, "espG1 snm5 Rv3866", "FUNCTION: Specific chaperone for cognate PE/PPE proteins. Plays an important role in preventing aggregation of PE/PPE dimers. {ECO:0000250|UniProtKB:O53943}.", 31),
    @("Rv1114", 3, "vapC32 Rv1114", "FUNCTION: Toxic component of a type II toxin-antitoxin (TA) system. An RNase (By similarity). Upon expression in M.smegmatis inhibits colony formation. Its toxic effect is neutralized by coexpression with cognate antitoxin VapB32. {ECO:0000255|HAMAP-Rule:MF_00265, ECO:0000269|PubMed:20011113}.", 31),
    @("Rv1382", 1, "Rv1382", "", 31),
    @("Rv0650", 1, "Rv0650", "", 31),
    @("Rv3829c", 1, "Rv3829c", "", 31),
    @("Rv0370c", 1, "Rv0370c", "", 31),
    @("Rv1141c", 1, "echA11 Rv1141c", "", 31)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
